# Add new mapping rows to the "mapping" sheet, then make the "data" sheet
# the active/selected sheet (instead of "mapping").

$wb = $excel.ActiveWorkbook
$mapping = $wb.Worksheets.Item("mapping")

# New rows: Column -> Location mapping entries.
$mapping.Range("A6").Value = "A"
$mapping.Range("B6").Value = "C9"

$mapping.Range("A7").Value = "E"
$mapping.Range("B7").Value = "C10"

$mapping.Range("A8").Value = "E"
$mapping.Range("B8").Value = "C11"

$mapping.Range("B8").Select()

# Make "data" the active sheet (tab selected).
$data = $wb.Worksheets.Item("data")
$data.Activate()
